# Word COM-interop script: Simplified Chinese -> Traditional Chinese
# translation update for "Email 10-1 [TEMPLATE] Partner email" (zh).
#
# Strategy: most segments are unique (or need the identical replacement at
# every occurrence) so a simple whole-document Find/Replace (wdReplaceAll)
# is used.  A couple of runs contain text that is ambiguous in isolation
# (the ", " that follows both [PARTNER NAME] and [CITY] needs two distinct
# replacements) - those are resolved positionally by first locating an
# unambiguous anchor (e.g. the bracketed placeholder) with Find and then
# acting on the Range that immediately follows it.

$d = $word.ActiveDocument

function ReplaceAll($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

# --- language switcher header (appears twice: hyperlink + plain label) ---
ReplaceAll "英语" "英文"
ReplaceAll " / 葡萄牙语 / 法语 / 泰语 / 越南语 / 西班牙语" " / 葡萄牙文 / 法文 / 泰文 / 越南文 / 西班牙文"

# --- info table (Introduction / Target audience) ---
ReplaceAll "简介" "簡介"
ReplaceAll "一封发送给参加活动的合作伙伴的邮件。 这封邮件将包含一个照片画廊，将通过 customer.io 发送。" "寄給參加活動的合作夥伴的電子郵件。 此電子郵件將包括照片畫廊，將通過 customer.io 發送。"
ReplaceAll "目标受众" "目標受眾"
ReplaceAll "活动参与者" "活動參加者"

# --- "Subject:" line ---
ReplaceAll "主题: " "Subject: "
ReplaceAll "感谢您参加 " "感謝您參加 "

# --- big heading ---
ReplaceAll "您使我们的活动圆满成功！ 🎉" "您使我們的活動圓滿成功！ 🎉"

# --- "Hi [PARTNER NAME], " greeting line ---
ReplaceAll "[合作伙伴姓名]" "[PARTNER NAME]"

# The literal ", " run exists twice (after the partner-name placeholder and
# after [CITY]); each needs a different replacement, so resolve positionally.
$docEnd = $d.Content.End
$afterPartner = $d.Range(0, $docEnd)
$foundPartner = $afterPartner.Find.Execute("PARTNER NAME", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
if ($foundPartner) {
    $commaAfterPartner = $d.Range($afterPartner.End, $docEnd)
    $foundComma1 = $commaAfterPartner.Find.Execute(", ", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
    if ($foundComma1) {
        $commaAfterPartner.Text = "， "
    }
}

# --- "Thank you for attending [EVENT NAME] in [CITY], [COUNTRY]." line ---
ReplaceAll " 于 " " 於 "

$docEnd = $d.Content.End
$afterCity = $d.Range(0, $docEnd)
$foundCity = $afterCity.Find.Execute("[CITY]", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
if ($foundCity) {
    $commaAfterCity = $d.Range($afterCity.End, $docEnd)
    $foundComma2 = $commaAfterCity.Find.Execute(", ", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
    if ($foundComma2) {
        $commaAfterCity.Text = "，"
    }
}

ReplaceAll "。 希望您玩得开心，很高兴认识您！" "。 希望您度過了愉快的時光，很高興認識您！"

# --- "browse the conference/workshop/alliance trip photos" paragraph ---
# (Split across 3 runs around a comment anchor; each run's own text changes.)
ReplaceAll "如需浏览 " "如需瀏覽會議/"
ReplaceAll "会议/研讨会/联盟之旅" "研討會/聯盟之旅"
ReplaceAll " 的照片和精彩片段，并随时了解我们举办的最新活动和计划，请关注我们：" "的照片和精彩片段，並隨時了解我們為您舉辦的最新活動和計劃，請關注我們的社交媒體帳戶："

# --- closing paragraph ---
ReplaceAll "希望这次活动能给您们带来和我们一样的启发，让我们继续共同成长！" "希望這次活動能給您們帶來和我們一樣的啟發，讓我們繼續共同成長！"

# --- comment body text ("choose one of these") ---
ReplaceAll "选择其中之一" "選擇其中一個"
